$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.980.15'
$ws.Cells.Item(2, 5).Value = '  +0.75%  '

$ws.Cells.Item(3, 4).Value = '2.631.06'
$ws.Cells.Item(3, 5).Value = '  +3.67%  '

$ws.Cells.Item(4, 5).Value = '  +0.23%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '517.66'
$c.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +2.01%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '145.24'
$c.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  +1.01%  '

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  +0.00%  '

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '0.567'
$c.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  +0.51%  '

$ws.Cells.Item(9, 4).Value = '2.656.80'
$ws.Cells.Item(9, 5).Value = '  +4.50%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '6.26'
$c.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  +2.76%  '

$ws.Cells.Item(11, 5).Value = '  +2.63%  '

$ws.Cells.Item(12, 5).Value = '  +1.74%  '

$ws.Cells.Item(13, 5).Value = '  -1.57%  '

$ws.Cells.Item(14, 4).Value = '3.119.15'
$ws.Cells.Item(14, 5).Value = '  +4.65%  '

$ws.Cells.Item(15, 4).Value = '58.993.23'
$ws.Cells.Item(15, 5).Value = '  +0.82%  '

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '20.99'
$c.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  +1.41%  '

$ws.Cells.Item(17, 5).Value = '  +1.49%  '

$ws.Cells.Item(18, 4).Value = '2.658.43'
$ws.Cells.Item(18, 5).Value = '  +4.70%  '

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '350.19'
$c.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  +3.88%  '

$ws.Cells.Item(20, 5).Value = '  +0.23%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '10.35'
$c.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +2.62%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '6.20'
$c.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +4.04%  '

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  +0.03%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '61.69'
$c.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +1.98%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '0.419'
$c.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +2.28%  '

$ws.Cells.Item(26, 4).Value = '2.763.04'
$ws.Cells.Item(26, 5).Value = '  +4.25%  '

$ws.Cells.Item(27, 2).Value = 'Kaspa'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '0.162'
$c.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +1.45%  '

$ws.Cells.Item(28, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '0.987'
$c.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  -1.20%  '

$ws.Cells.Item(29, 5).Value = '  +2.48%  '

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '7.14'
$c.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +2.70%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -0.01%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '6.27'
$c.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  +7.42%  '

$ws.Cells.Item(33, 5).Value = '  +2.64%  '

$ws.Cells.Item(34, 5).Value = '  +2.80%  '

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '149.83'
$c.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +0.08%  '

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '0.968'
$c.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  +6.08%  '

$ws.Cells.Item(37, 5).Value = '  +3.17%  '

$ws.Cells.Item(38, 5).Value = '  +2.54%  '

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '36.77'
$c.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +1.92%  '

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '3.71'
$c.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +5.17%  '

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '1.42'
$c.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  +1.72%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '277.68'
$c.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  -2.19%  '

$ws.Cells.Item(44, 5).Value = '  +0.02%  '

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.610'
$c.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  +1.75%  '

$ws.Cells.Item(46, 5).Value = '  -1.18%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '19.56'
$c.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  +4.64%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '0.0525'
$c.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -1.08%  '

$ws.Cells.Item(49, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '10.29'
$c.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -0.01%  '

$ws.Cells.Item(50, 2).Value = 'VeChain'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.0230'
$c.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  +1.38%  '

$ws.Cells.Item(51, 4).Value = '1.996.22'
$ws.Cells.Item(51, 5).Value = '  +4.75%  '
